$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.287107
$ws.Range("H2").Value = 0.861321
$ws.Range("M2").Value = 1.345575333333333
$ws.Range("N2").Value = 4.036726
$ws.Range("O2").Value = 0.05720258880468643
$ws.Range("P2").Value = 0.05720258880468642
$ws.Range("Q2").Value = 0.3863240972273334
$ws.Range("R2").Value = 3.476916875046
$ws.Range("S2").Value = 0.05720258880468643
$ws.Range("T2").Value = 0.05720258880468642

# Row 3
$ws.Range("G3").Value = 0.287107
$ws.Range("H3").Value = 0.861321
$ws.Range("O3").Value = 0.7063377274242425
$ws.Range("P3").Value = 0.7063377274242425
$ws.Range("Q3").Value = 4.770331038976
$ws.Range("R3").Value = 42.932979350784
$ws.Range("S3").Value = 0.7063377274242425
$ws.Range("T3").Value = 0.7063377274242425

# Row 4
$ws.Range("G4").Value = 0.287107
$ws.Range("H4").Value = 0.861321
$ws.Range("M4").Value = 0.472618
$ws.Range("N4").Value = 1.417854
$ws.Range("O4").Value = 0.02009175736650936
$ws.Range("P4").Value = 0.02009175736650936
$ws.Range("Q4").Value = 0.135691936126
$ws.Range("R4").Value = 1.221227425134
$ws.Range("S4").Value = 0.02009175736650936
$ws.Range("T4").Value = 0.02009175736650936

# Row 5
$ws.Range("G5").Value = 0.287107
$ws.Range("H5").Value = 0.861321
$ws.Range("M5").Value = 1.403935
$ws.Range("N5").Value = 4.211805
$ws.Range("O5").Value = 0.05968355284468708
$ws.Range("P5").Value = 0.05968355284468708
$ws.Range("Q5").Value = 0.403079566045
$ws.Range("R5").Value = 3.627716094405
$ws.Range("S5").Value = 0.05968355284468708
$ws.Range("T5").Value = 0.05968355284468708

# Row 6
$ws.Range("G6").Value = 0.287107
$ws.Range("H6").Value = 0.861321
$ws.Range("M6").Value = 2.832734666666667
$ws.Range("N6").Value = 8.498204000000001
$ws.Range("O6").Value = 0.1204241429788253
$ws.Range("P6").Value = 0.1204241429788253
$ws.Range("Q6").Value = 0.8132979519426667
$ws.Range("R6").Value = 7.319681567484001
$ws.Range("S6").Value = 0.1204241429788253
$ws.Range("T6").Value = 0.1204241429788253

# Row 7
$ws.Range("G7").Value = 0.287107
$ws.Range("H7").Value = 0.861321
$ws.Range("M7").Value = 0.8529486666666667
$ws.Range("N7").Value = 2.558846
$ws.Range("O7").Value = 0.03626023058104926
$ws.Range("P7").Value = 0.03626023058104925
$ws.Range("Q7").Value = 0.2448875328406667
$ws.Range("R7").Value = 2.203987795566
$ws.Range("S7").Value = 0.03626023058104926
$ws.Range("T7").Value = 0.03626023058104925
